# Update metadata values on the "Metadata" sheet of the CodeSystem workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# URL: https://hl7.fr/fhir/fr/medication/... -> https://hl7.fr/ig/fhir/medication/...
$ws.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-medication-reconciliation-discrepancy"

# Name: FrMedicationReconciliationDiscrepancy -> FRMedicationReconciliationDiscrepancy
$ws.Range("B4").Value = "FRMedicationReconciliationDiscrepancy"

# Date: 2025-04-10T15:35:36+00:00 -> 2026-01-15T08:54:26+00:00
$ws.Range("B8").Value = "2026-01-15T08:54:26+00:00"

# Jurisdiction: (empty) -> FRANCE
$ws.Range("B11").Value = "FRANCE"
